# Generate Report for Handoff
#
# The "abb6b744-1871-4545-a88f-6fac0f1d6f99" localization record moved from
# "Handed back: in sync with en-US" to "Ready for handoff": the Overview
# sheet's summary status, and each language sheet's detail row, are updated
# with the new status, refreshed handoff timestamps, and (for this item) a
# stale-handback-version error message. The Error Detail column is also
# widened to fit the new message.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for abb6b744-...md (row 3) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-26 12:48:50"

# --- zh-cn sheet: row for abb6b744-...md (row 3) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("H3").Value = "2016-08-26 12:48:46"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/900753929cc2e35a0fd4b6eda84777c8d29cad6c/e2e/abb6b744-1871-4545-a88f-6fac0f1d6f99.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a907b4f022d7dce7ec55aaba4a9e2961971e643/e2e/abb6b744-1871-4545-a88f-6fac0f1d6f99.md."
# Widen the Error Detail column (P) to fit the new message (matches width=40)
$ws.Range("P1").ColumnWidth = 39.1667

# --- de-de sheet: row for abb6b744-...md (row 3) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("H3").Value = "2016-08-26 12:48:50"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/900753929cc2e35a0fd4b6eda84777c8d29cad6c/e2e/abb6b744-1871-4545-a88f-6fac0f1d6f99.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a907b4f022d7dce7ec55aaba4a9e2961971e643/e2e/abb6b744-1871-4545-a88f-6fac0f1d6f99.md."
# Widen the Error Detail column (P) to fit the new message (matches width=40)
$ws.Range("P1").ColumnWidth = 39.1667
